$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 981.9666999999999
$ws.Range("I40").Value = 959.96155
$ws.Range("J40").Value = 1125
$ws.Range("K40").Value = 959.96155
$ws.Range("L40").Value = 1125
$ws.Range("M40").Value = -784.96155
$ws.Range("N40").Value = -1475
$ws.Range("H129").Value = 747.7143
$ws.Range("I129").Value = 686
$ws.Range("J129").Value = 830
$ws.Range("K129").Value = 2058
$ws.Range("L129").Value = 2490
$ws.Range("M129").Value = 2942
$ws.Range("N129").Value = -12490

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9594.777
$ws.Range("I32").Value = 6289.113
$ws.Range("J32").Value = 30089.9
$ws.Range("K32").Value = 6289.113
$ws.Range("L32").Value = 30089.9
$ws.Range("M32").Value = -6002.113
$ws.Range("N32").Value = -30663.9
$ws.Range("H74").Value = 13047481
$ws.Range("I74").Value = 17648010
$ws.Range("J74").Value = 12647.417
$ws.Range("K74").Value = 17648010
$ws.Range("L74").Value = 12647.417
$ws.Range("M74").Value = -17647136
$ws.Range("N74").Value = -14395.417
$ws.Range("H77").Value = 13047481
$ws.Range("I77").Value = 17648010
$ws.Range("J77").Value = 12647.417
$ws.Range("K77").Value = 88240050
$ws.Range("L77").Value = 63237.085
$ws.Range("M77").Value = -88235682
$ws.Range("N77").Value = -71973.08499999999
$ws.Range("H102").Value = 2119.9285
$ws.Range("I102").Value = 1888.7778
$ws.Range("J102").Value = 2536
$ws.Range("K102").Value = 1888.7778
$ws.Range("L102").Value = 2536
$ws.Range("M102").Value = -266.7778000000001
$ws.Range("N102").Value = -5780
$ws.Range("H122").Value = 1358.8214
$ws.Range("I122").Value = 997.5238000000001
$ws.Range("J122").Value = 2442.7144
$ws.Range("K122").Value = 2992.5714
$ws.Range("L122").Value = 7328.1432
$ws.Range("M122").Value = -542.5714000000003
$ws.Range("N122").Value = -12228.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1536
$ws.Range("I99").Value = 1300
$ws.Range("J99").Value = 1595
$ws.Range("K99").Value = 1300
$ws.Range("L99").Value = 1595
$ws.Range("M99").Value = 198
$ws.Range("N99").Value = -4591
$ws.Range("H134").Value = 9141667
$ws.Range("I134").Value = 9576937
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 28730811
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -28728276
$ws.Range("N134").Value = -8070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1558.6333
$ws.Range("I31").Value = 1083.2273
$ws.Range("J31").Value = 2866
$ws.Range("K31").Value = 1083.2273
$ws.Range("L31").Value = 2866
$ws.Range("M31").Value = -788.2273
$ws.Range("N31").Value = -3456
$ws.Range("H34").Value = 1558.6333
$ws.Range("I34").Value = 1083.2273
$ws.Range("J34").Value = 2866
$ws.Range("K34").Value = 1083.2273
$ws.Range("L34").Value = 2866
$ws.Range("M34").Value = -881.2273
$ws.Range("N34").Value = -3270
$ws.Range("H99").Value = 1577.8422
$ws.Range("I99").Value = 1465
$ws.Range("J99").Value = 1771.2858
$ws.Range("K99").Value = 1465
$ws.Range("L99").Value = 1771.2858
$ws.Range("M99").Value = 33
$ws.Range("N99").Value = -4767.2858
$ws.Range("H126").Value = 1577.8422
$ws.Range("I126").Value = 1465
$ws.Range("J126").Value = 1771.2858
$ws.Range("K126").Value = 4395
$ws.Range("L126").Value = 5313.857400000001
$ws.Range("M126").Value = -1925
$ws.Range("N126").Value = -10253.8574
$ws.Range("H132").Value = 2729.0476
$ws.Range("I132").Value = 2726.5334
$ws.Range("J132").Value = 2735.3333
$ws.Range("K132").Value = 8179.600199999999
$ws.Range("L132").Value = 8205.999899999999
$ws.Range("M132").Value = -5649.600199999999
$ws.Range("N132").Value = -13265.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1928.179
$ws.Range("I131").Value = 5550.8
$ws.Range("J131").Value = 1501.9883
$ws.Range("K131").Value = 16652.4
$ws.Range("L131").Value = 4505.9649
$ws.Range("M131").Value = -11612.4
$ws.Range("N131").Value = -14585.9649
$ws.Range("H134").Value = 4865.7144
$ws.Range("I134").Value = 3811.6667
$ws.Range("J134").Value = 5656.25
$ws.Range("K134").Value = 11435.0001
$ws.Range("L134").Value = 16968.75
$ws.Range("M134").Value = -6365.000100000001
$ws.Range("N134").Value = -27108.75
$ws.Range("H137").Value = 55099.6
$ws.Range("I137").Value = 4286.6665
$ws.Range("J137").Value = 64066.59
$ws.Range("K137").Value = 12859.9995
$ws.Range("L137").Value = 192199.77
$ws.Range("M137").Value = -7759.999500000002
$ws.Range("N137").Value = -202399.77

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 725.8261
$ws.Range("I107").Value = 587
$ws.Range("J107").Value = 853.0833
$ws.Range("K107").Value = 587
$ws.Range("L107").Value = 853.0833
$ws.Range("M107").Value = 1333
$ws.Range("N107").Value = -4693.0833
$ws.Range("H126").Value = 2571.1428
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2799.6
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 8398.799999999999
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -13338.8
$ws.Range("H132").Value = 1793.2258
$ws.Range("I132").Value = 1550.12
$ws.Range("J132").Value = 2806.1667
$ws.Range("K132").Value = 4650.36
$ws.Range("L132").Value = 8418.500100000001
$ws.Range("M132").Value = -2120.36
$ws.Range("N132").Value = -13478.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1010.7143
$ws.Range("I61").Value = 875
$ws.Range("J61").Value = 1191.6666
$ws.Range("K61").Value = 875
$ws.Range("L61").Value = 1191.6666
$ws.Range("M61").Value = -673
$ws.Range("N61").Value = -1595.6666
$ws.Range("H93").Value = 1644.1818
$ws.Range("I93").Value = 1528
$ws.Range("J93").Value = 1847.5
$ws.Range("K93").Value = 1528
$ws.Range("L93").Value = 1847.5
$ws.Range("M93").Value = -280
$ws.Range("N93").Value = -4343.5
$ws.Range("H100").Value = 1848.8334
$ws.Range("I100").Value = 1750.25
$ws.Range("J100").Value = 2046
$ws.Range("K100").Value = 1750.25
$ws.Range("L100").Value = 2046
$ws.Range("M100").Value = -1209.25
$ws.Range("N100").Value = -3128
$ws.Range("H113").Value = 1010.7143
$ws.Range("I113").Value = 875
$ws.Range("J113").Value = 1191.6666
$ws.Range("K113").Value = 875
$ws.Range("L113").Value = 1191.6666
$ws.Range("M113").Value = 1295
$ws.Range("N113").Value = -5531.6666
$ws.Range("H122").Value = 2533.2222
$ws.Range("I122").Value = 2159.8
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6479.400000000001
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4029.400000000001
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 6148.1113
$ws.Range("I132").Value = 6312.5415
$ws.Range("J132").Value = 4832.6665
$ws.Range("K132").Value = 18937.6245
$ws.Range("L132").Value = 14497.9995
$ws.Range("M132").Value = -16407.6245
$ws.Range("N132").Value = -19557.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4530.061
$ws.Range("I132").Value = 4865.952
$ws.Range("J132").Value = 2514.7144
$ws.Range("K132").Value = 14597.856
$ws.Range("L132").Value = 7544.1432
$ws.Range("M132").Value = -12067.856
$ws.Range("N132").Value = -12604.1432
$ws.Range("H136").Value = 6686.4
$ws.Range("I136").Value = 8119.222
$ws.Range("J136").Value = 1850.625
$ws.Range("K136").Value = 24357.666
$ws.Range("L136").Value = 5551.875
$ws.Range("M136").Value = -21807.666
$ws.Range("N136").Value = -10651.875
